$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.882.38'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '2.101.41'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.82'
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.653'
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.00'
$ws.Range("E8").Value = '  -4.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.77'
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.370'
$ws.Range("E10").Value = '  -4.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0773'
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.15'
$ws.Range("E13").Value = '  -5.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.890'
$ws.Range("E14").Value = '  +6.08%  '
$ws.Range("D15").Value = '2.394.51'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.60'
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("D17").Value = '2.052.18'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '36.846.52'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.56'
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.57'
$ws.Range("E20").Value = '  -2.23%  '
$ws.Range("D21").Value = '0.0₃0883'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.52'
$ws.Range("E22").Value = '  +2.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.54'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.91'
$ws.Range("E26").Value = '  +5.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.19'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.73'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '21.03'
$ws.Range("E29").Value = '  +4.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.40'
$ws.Range("E30").Value = '  +11.39%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.20'
$ws.Range("E32").Value = '  +6.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.76'
$ws.Range("E33").Value = '  +4.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0614'
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.45'
$ws.Range("E35").Value = '  +6.82%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  +4.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0848'
$ws.Range("E38").Value = '  -5.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.29'
$ws.Range("E39").Value = '  -3.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.18'
$ws.Range("E40").Value = '  +2.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.93'
$ws.Range("E41").Value = '  -5.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0222'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0964'
$ws.Range("E43").Value = '  -7.22%  '
$ws.Range("E44").Value = '  -8.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '97.17'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.46'
$ws.Range("E46").Value = '  -5.30%  '
$ws.Range("D47").Value = '1.359.09'
$ws.Range("E47").Value = '  +5.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.43'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E49").Value = '  +3.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.90'
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").Value = '2.283.23'
$ws.Range("E51").Value = '  +1.74%  '
